# Paket 3 selesai dikoreksi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new score values for "Paket 3" columns (H, I, J) for several rows.
$ws.Range("I2").Value = 0

$ws.Range("H3").Value = 7
$ws.Range("I3").Value = 6

$ws.Range("H4").Value = 7
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 7

$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 0

$ws.Range("I8").Value = 3

$ws.Range("H9").Value = 7

# Update the active selection to reflect where the user finished editing.
$ws.Range("J5").Select()
